# Update the regression-results table (Crisis and Credit Allocation
# coefficients) in place. Values are stored as text in the sheet (some
# carry significance stars, e.g. "0.46***"), so cells whose new value
# would otherwise look like a plain number ("0.17", "-0.01", "-0.09",
# "0.98", "-0.89") must be forced to Text format before assignment,
# otherwise Excel auto-converts them to numeric cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$numericLooking = @("B2", "B3", "B4", "C4", "D2")
foreach ($addr in $numericLooking) {
    $ws.Range($addr).NumberFormat = "@"
}

# Assign column-by-column (B then C then D) to mirror the original
# layout/order of the underlying shared-string table.
$ws.Range("B2").Value = "0.17"
$ws.Range("B3").Value = "-0.01"
$ws.Range("B4").Value = "-0.09"

$ws.Range("C2").Value = "44.29***"
$ws.Range("C3").Value = "2.21***"
$ws.Range("C4").Value = "0.98"

$ws.Range("D2").Value = "-0.89"
$ws.Range("D3").Value = "0.46***"
$ws.Range("D4").Value = "0.82*"
